$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text formatting (avoid Excel auto-
# converting numeric-looking strings like "30.099.40" or "0.5162" into
# numbers) by forcing a Text number format before assigning the value.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.099.40'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.107.50'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '348.33'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5162'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4449'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -2.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.63'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -4.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08977'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.50%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '25.73'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.118.02'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.313'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.749'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '99.13'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001151'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.92'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +8.05%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.06688'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.251'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.210.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.85'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.339'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.69%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.362.99'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.543'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.64'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.83'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.178'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.81%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.24%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.637'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.262'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.962'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.35'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.922'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02580'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06830'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2312'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.66'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.71%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6829'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.284'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.33'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.309'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6389'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000366'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.74%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.654'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.92%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '82.84'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07230'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.53%  '
